$d = $word.ActiveDocument

# The document has one section whose header/footer pictures carry
# "name" metadata (wp:docPr / pic:cNvPr) that doesn't match the asset
# they actually are:
#   - the Pearson logo (footers) is labelled "image1.png" -> should be "image2.png"
#   - the BTec logo    (headers) is labelled "image2.jpg" -> should be "image1.jpg"
# Relabel every occurrence (primary + first-page header/footer) without
# touching anything else about the pictures (size, alt text/description,
# embed relationship, etc. all stay the same).

$sec = $d.Sections.Item(1)

# --- Footers: Pearson logo, "image1.png" -> "image2.png" ---
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $rng = $ftr.Range
        for ($i = 1; $i -le $rng.InlineShapes.Count; $i++) {
            $rng.InlineShapes.Item($i).Name = "image2.png"
        }
    }
}

# --- Headers: BTec logo, "image2.jpg" -> "image1.jpg" ---
for ($h = 1; $h -le $sec.Headers.Count; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $rng = $hdr.Range
        for ($i = 1; $i -le $rng.InlineShapes.Count; $i++) {
            $rng.InlineShapes.Item($i).Name = "image1.jpg"
        }
    }
}
